$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1398
$ws.Range("C3").Value = 3421
$ws.Range("D3").Value = 4971
$ws.Range("E3").Value = 28700
$ws.Range("F3").Value = 31600
$ws.Range("G3").Value = 30000
$ws.Range("B4").Value = 7062
$ws.Range("C4").Value = 13900
$ws.Range("D4").Value = 22300
$ws.Range("E4").Value = 25800
$ws.Range("F4").Value = 32500
$ws.Range("G4").Value = 29300
$ws.Range("B5").Value = 7198
$ws.Range("C5").Value = 14400
$ws.Range("D5").Value = 21700
$ws.Range("E5").Value = 23700
$ws.Range("F5").Value = 26100
$ws.Range("G5").Value = 25900
$ws.Range("B10").Value = 19600
$ws.Range("C10").Value = 30400
$ws.Range("D10").Value = 56200
$ws.Range("E10").Value = 107000
$ws.Range("F10").Value = 187000
$ws.Range("G10").Value = 286000
$ws.Range("B11").Value = 13300
$ws.Range("C11").Value = 26400
$ws.Range("D11").Value = 52700
$ws.Range("E11").Value = 140000
$ws.Range("F11").Value = 192000
$ws.Range("G11").Value = 181000
$ws.Range("B12").Value = 15600
$ws.Range("C12").Value = 41000
$ws.Range("D12").Value = 81600
$ws.Range("E12").Value = 144000
$ws.Range("F12").Value = 180000
$ws.Range("G12").Value = 180000
$ws.Range("B17").Value = 4751
$ws.Range("C17").Value = 5538
$ws.Range("D17").Value = 5734
$ws.Range("E17").Value = 5993
$ws.Range("F17").Value = 5514
$ws.Range("G17").Value = 4914
$ws.Range("B18").Value = 4011
$ws.Range("C18").Value = 5067
$ws.Range("D18").Value = 5506
$ws.Range("E18").Value = 5817
$ws.Range("F18").Value = 5673
$ws.Range("G18").Value = 4843
$ws.Range("B19").Value = 4015
$ws.Range("C19").Value = 5128
$ws.Range("D19").Value = 5496
$ws.Range("E19").Value = 5955
$ws.Range("F19").Value = 6089
$ws.Range("G19").Value = 4773
$ws.Range("B24").Value = 381000
$ws.Range("C24").Value = 484000
$ws.Range("D24").Value = 622000
$ws.Range("E24").Value = 647000
$ws.Range("F24").Value = 200000
$ws.Range("G24").Value = 122000
$ws.Range("B25").Value = 117000
$ws.Range("C25").Value = 201000
$ws.Range("D25").Value = 294000
$ws.Range("E25").Value = 175000
$ws.Range("F25").Value = 129000
$ws.Range("G25").Value = 130000
$ws.Range("B26").Value = 152000
$ws.Range("C26").Value = 265000
$ws.Range("D26").Value = 368000
$ws.Range("E26").Value = 125000
$ws.Range("F26").Value = 139000
$ws.Range("G26").Value = 129000
$ws.Range("B31").Value = 3312
$ws.Range("C31").Value = 5221
$ws.Range("D31").Value = 12300
$ws.Range("E31").Value = 7488
$ws.Range("F31").Value = 10500
$ws.Range("G31").Value = 16200
$ws.Range("B32").Value = 3343
$ws.Range("C32").Value = 4895
$ws.Range("D32").Value = 8696
$ws.Range("E32").Value = 11400
$ws.Range("F32").Value = 10000
$ws.Range("G32").Value = 14900
$ws.Range("B33").Value = 3893
$ws.Range("C33").Value = 4641
$ws.Range("D33").Value = 12400
$ws.Range("E33").Value = 7434
$ws.Range("F33").Value = 11300
$ws.Range("G33").Value = 15200
$ws.Range("B38").Value = 127000
$ws.Range("C38").Value = 161000
$ws.Range("D38").Value = 411000
$ws.Range("E38").Value = 266000
$ws.Range("F38").Value = 343000
$ws.Range("G38").Value = 463000
$ws.Range("B39").Value = 112000
$ws.Range("C39").Value = 135000
$ws.Range("D39").Value = 245000
$ws.Range("E39").Value = 344000
$ws.Range("F39").Value = 353000
$ws.Range("G39").Value = 502000
$ws.Range("B40").Value = 144000
$ws.Range("C40").Value = 188000
$ws.Range("D40").Value = 394000
$ws.Range("E40").Value = 240000
$ws.Range("F40").Value = 689000
$ws.Range("G40").Value = 468000
$ws.Range("B45").Value = 10800
$ws.Range("C45").Value = 12600
$ws.Range("D45").Value = 13000
$ws.Range("E45").Value = 15000
$ws.Range("F45").Value = 15200
$ws.Range("G45").Value = 3030
$ws.Range("B46").Value = 2861
$ws.Range("C46").Value = 3155
$ws.Range("D46").Value = 3133
$ws.Range("E46").Value = 3253
$ws.Range("F46").Value = 3392
$ws.Range("G46").Value = 3663
$ws.Range("B47").Value = 2890
$ws.Range("C47").Value = 3130
$ws.Range("D47").Value = 3459
$ws.Range("E47").Value = 3750
$ws.Range("F47").Value = 3867
$ws.Range("G47").Value = 3880
$ws.Range("B52").Value = 322000
$ws.Range("C52").Value = 386000
$ws.Range("D52").Value = 450000
$ws.Range("E52").Value = 476000
$ws.Range("F52").Value = 490000
$ws.Range("G52").Value = 492000
$ws.Range("B53").Value = 99900
$ws.Range("C53").Value = 188000
$ws.Range("D53").Value = 278000
$ws.Range("E53").Value = 401000
$ws.Range("F53").Value = 456000
$ws.Range("G53").Value = 91300
$ws.Range("B54").Value = 58900
$ws.Range("C54").Value = 73500
$ws.Range("D54").Value = 78500
$ws.Range("E54").Value = 85800
$ws.Range("F54").Value = 99900
$ws.Range("G54").Value = 112000
